$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (sharedStrings rich-text runs) ---
# A8: "Volume 31   Number  39" -> "...40"
$ws.Range("A8").Characters(21, 2).Text = "40"
# C9: "Report Covering the Week  9/23/2024  Through  9/29/2024"
#     -> "...9/30/2024  Through  10/6/2024"
$ws.Range("C9").Characters(27, 9).Text = "9/30/2024"
$ws.Range("C9").Characters(47, 9).Text = "10/6/2024"

# --- Precinct crime-stat table updates (rows 15-28) ---
# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("I14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 11.111111111111
$ws.Range("L15").Value = -37.5
$ws.Range("N15").Value = -52.380952380952
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 175
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 75
$ws.Range("L16").Value = 47.058823529411
$ws.Range("M16").Value = 75
$ws.Range("N16").Value = -83.269598470363
# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 231
$ws.Range("J17").Value = 139
$ws.Range("K17").Value = 66.187050359712
$ws.Range("L17").Value = 54
$ws.Range("M17").Value = 115.88785046729
$ws.Range("N17").Value = -39.050131926121
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 27.272727272727
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 54.545454545454
$ws.Range("L18").Value = -14.044943820224
$ws.Range("M18").Value = 5.517241379310
$ws.Range("N18").Value = -90.721649484536
# Row 19
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = -11.428571428571
$ws.Range("F19").Value = 115
$ws.Range("G19").Value = 141
$ws.Range("H19").Value = -18.439716312056
$ws.Range("I19").Value = 1364
$ws.Range("J19").Value = 1385
$ws.Range("K19").Value = -1.516245487364
$ws.Range("L19").Value = -7.588075880758
$ws.Range("M19").Value = 8.685258964143
$ws.Range("N19").Value = -76.182992840928
# Row 20
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -30
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -50.537634408602
$ws.Range("N20").Value = -88.082901554404
# Row 21
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 2.083333333333
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = -10.396039603960
$ws.Range("I21").Value = 1982
$ws.Range("J21").Value = 1803
$ws.Range("K21").Value = 9.927897947864
$ws.Range("L21").Value = -2.508607968519
$ws.Range("M21").Value = 19.253910950661
$ws.Range("N21").Value = -78.498589715773
# Row 22
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -83.333333333333
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 10
$ws.Range("I22").Value = 60
$ws.Range("J22").Value = 59
$ws.Range("K22").Value = 1.694915254237
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 30.434782608695
# Row 24
$ws.Range("C24").Value = 71
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = 22.413793103448
$ws.Range("F24").Value = 281
$ws.Range("G24").Value = 226
$ws.Range("H24").Value = 24.336283185840
$ws.Range("I24").Value = 2443
$ws.Range("J24").Value = 2094
$ws.Range("K24").Value = 16.666666666666
$ws.Range("L24").Value = 20.522940305870
$ws.Range("M24").Value = 66.41689373297
# Row 25
$ws.Range("C25").Value = 62
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = 26.530612244898
$ws.Range("F25").Value = 240
$ws.Range("G25").Value = 205
$ws.Range("H25").Value = 17.073170731707
$ws.Range("I25").Value = 2319
$ws.Range("J25").Value = 2054
$ws.Range("K25").Value = 12.901655306718
$ws.Range("L25").Value = 6.964944649446
# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 65
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 578
$ws.Range("J26").Value = 529
$ws.Range("K26").Value = 9.262759924385
$ws.Range("L26").Value = 32.568807339449
$ws.Range("M26").Value = 67.052023121387
# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -10
$ws.Range("L27").Value = -35.714285714285
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 84
$ws.Range("J28").Value = 78
$ws.Range("K28").Value = 7.692307692307
$ws.Range("L28").Value = 16.666666666666
